$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the text of the closing italic paragraph (old meta-description
#    sentence) with the new image-prompt text, keeping its formatting intact.
#    This is done FIRST, while the sentence still occurs only once in the
#    document, so the replace cannot accidentally hit any other location.
# ---------------------------------------------------------------------------
$oldSentence = "Read our review of the Asgardian Stones slot game featuring gameplay, bonus features, graphics and sound, and betting options. Play for free today."
$newSentence = "Create an eye-catching feature image for Asgardian Stones that features a happy Maya warrior with glasses. The image should be in cartoon style and should convey the excitement and thrill of the game. You could include elements from Norse mythology, such as the Asgardian Stones or symbols of power and strength to add to the theme. Make the image bright and colorful with bold outlines to make it pop. The Maya warrior should be depicted as having fun and enjoying the game to encourage players to give it a try. The image should be of high quality and clearly convey the message that this is a fun and exciting game to play."

$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the paragraph that used to hold the bold title text near the
#    end of the document ("Play Asgardian Stones for Free - Review and
#    Features"). Deleting the paragraph's Range (without extending it)
#    removes the text together with its paragraph mark, merging the
#    following paragraph upward.
# ---------------------------------------------------------------------------
$target = "Play Asgardian Stones for Free - Review and Features"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $target) {
        $p.Range.Delete() | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph "Play Asgardian Stones for Free - Review and
#    Features". The new paragraph has no paragraph style, and contains an
#    empty run, a bold run "Meta description", and a regular run with the
#    rest of the sentence.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstRange = $firstPara.Range
$firstRange.Collapse(0)
$firstRange.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range

$metaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r/>
<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
<w:r><w:t>: Read our review of the Asgardian Stones slot game featuring gameplay, bonus features, graphics and sound, and betting options. Play for free today.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$metaRange.InsertXML($metaXml) | Out-Null
